# GAF Stundenplan löschen erfasst
# Adds a new row (row 5) to the "Geschäftsanwendungfall" sheet describing
# the "Aufruf der Funktion Stundenplan löschen" use case that documents
# calling the delete-schedule function.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row content -------------------------------------------------
$ws.Range("A5").Value = "BUC4"
$ws.Range("B5").Value = "Stundenplan löschen"
$ws.Range("C5").Value = "Aufruf der Funktion Stundenplan löschen"
$ws.Range("D5").Value = " - Eine Planung wurde bereits durchgeführt.`n - Eine neue Planung soll durchgeführt werden."
$ws.Range("E5").Value = " - Mitarbeiter Verwaltung"
$ws.Range("G5").Value = "1) Benutzer versucht eine neue Stundenplanberechnung durchzuführen.`n2) Eine Hinweismeldung erscheitn, dass die alte Planung zuerst gelöscht werden muss.`n3) Über den Button `"Vorhandenen Plan löschen`" löscht das System die bestehende Planung."
$ws.Range("I5").Value = "Der bestehende Stundenplan wurde aus der Datenbank entfernt."

# --- Formatting: mirror the look of the row above (row 4) ------------
# Row 4 uses border + vertical-center on every cell, with wrap text only
# on the "paragraph" style columns (C, D, G, I).
$row5 = $ws.Range("A5:I5")
$row5.Borders().LineStyle = 1
$row5.VerticalAlignment = -4108

$ws.Range("C5").WrapText = $true
$ws.Range("D5").WrapText = $true
$ws.Range("G5").WrapText = $true
$ws.Range("I5").WrapText = $true

$ws.Rows.Item(5).RowHeight = 135

# --- View state: match selection reported after the edit -------------
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 3
$ws.Range("I6").Select()
